$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

$ws.Cells.Item($row, 1).Value = 45968
$ws.Cells.Item($row, 2).Value = "22,1388"
$ws.Cells.Item($row, 3).Value = "15,8667"
$ws.Cells.Item($row, 4).Value = "15,2818"
$ws.Cells.Item($row, 5).Value = "15,2818"

$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
